# Fixed naive component forecaster bug - Presentation state 11.02.
#
# Each data row (2..16) holds a rolling window of QoQ errors in columns
# B:K. A new error value is produced each period and prepended into
# column B, pushing all older values one column to the right; the
# oldest value that would fall past column K is dropped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New error value to prepend into column B for each data row (2..16).
$newValues = @{
    2  = -0.7322633397437844
    3  = 0.2703549766394939
    4  = -1.355327161308811
    5  = 1.651602845777944
    6  = 0.3282974736644749
    7  = 0.7356582956163805
    8  = 0.1181882633125878
    9  = 0.7543890506736601
    10 = -0.1543252035281459
    11 = 0.2293445564577608
    12 = 0.2201546830999171
    13 = 0.314534851581486
    14 = -0.5970339283829468
    15 = 0.1550649743121164
    16 = -0.1624199859130616
}

$firstDataCol = 2   # column B
$lastCol = 11       # column K

for ($r = 2; $r -le 16; $r++) {

    # Find the last populated column in B:K for this row before shifting.
    $lastUsed = $firstDataCol - 1
    for ($c = $firstDataCol; $c -le $lastCol; $c++) {
        $cellVal = $ws.Cells.Item($r, $c).Value2
        if ($cellVal -ne $null) {
            $lastUsed = $c
        }
    }

    # Shift existing values one column to the right, starting from the
    # rightmost populated column and working left so values aren't
    # clobbered before they are read. Any value already in the last
    # column (K) is simply overwritten/dropped.
    for ($c = $lastUsed; $c -ge $firstDataCol; $c--) {
        $srcVal = $ws.Cells.Item($r, $c).Value2
        $destCol = $c + 1
        if ($destCol -le $lastCol) {
            $ws.Cells.Item($r, $destCol).Value2 = $srcVal
        }
    }

    # Prepend the new value into column B.
    $ws.Cells.Item($r, $firstDataCol).Value2 = $newValues[$r]
}
